# Apply the "0.0.5 with -99xx codes" revision to the
# controlled-terminology sheet of the mvreason workbook.
#
# Column A (identifier) and column E (broader-term lookup) are driven by
# formulas that key off column B, so updating B21:B45 (and B11/B12) is
# sufficient -- Excel recalculates A/E automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bump the version / last-updated metadata (B11:B12) ---------------
$ws.Range("B11").Value = "2022-03-23T15:01:00-08:00"
$ws.Range("B12").Value = "0.0.5"

# --- append the RADx "-99xx" numeric codes to each reason label -------
$labels = @{
    21 = "Reason Unknown (-9999)"
    22 = "Not Sent to Data Hub (-9980)"
    23 = "Data Transfer Agreement (-9981)"
    24 = "No Participant Consent to Share (-9982)"
    25 = "Not Available Or Mappable (-9983)"
    26 = "Data Lost Or Inaccessible (-9984)"
    27 = "Data Invalid (-9985)"
    28 = "Anonymization Or Privacy Concerns (-9986)"
    29 = "Other Unsent Reason Not Specified (-9987)"
    30 = "Not Entered By Originator (-9960)"
    31 = "Omitted This Value (-9961)"
    32 = "Originator Chose to Omit (-9962)"
    33 = "Question Not Applicable (-9963)"
    34 = "Answer Not Known (-9964)"
    35 = "Record Not Provided (-9965)"
    36 = "All Originators Omitted Element (-9966)"
    37 = "CDE Omitted with Exception (-9967)"
    38 = "Other Unentered Reason Not Specified (-9968)"
    39 = "Not Presented to Participant (-9940)"
    40 = "Skip Logic (-9941)"
    41 = "No Participant Consent to Ask (-9942)"
    42 = "CDE Not Presented Due to Exception (-9943)"
    43 = "Element Never Presented for Collection (-9944)"
    44 = "Process Error (-9945)"
    45 = "Other Unpresented Reason Not Specified (-9946)"
}

foreach ($row in $labels.Keys) {
    $ws.Range("B$row").Value = $labels[$row]
}

# --- move the saved cell selection to match the author's last edit ----
$ws.Range("B11").Select()

$wb.Save()
